$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 34 ("TimeSpan to Text:" / percentage test row),
# pushing the existing "@ format to Number:" row (and everything below it)
# down by one.
$ws.Rows.Item(34).Insert()

# Fill in the new row: label + a percentage value that must be parsed as a
# number when the cell's data type is set to numeric (fixes #687).
$ws.Range("B34").Value = "Percentage Text to Number:"
$ws.Range("C34").Value = 0.5512
$ws.Range("C34").NumberFormat = "0.00%"

# The row label column was widened slightly to fit the new, longer caption.
$ws.Columns.Item(2).ColumnWidth = 25.1
